$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.969.46"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").Value = "'2.586.84"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'521.45"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "'139.65"
$ws.Range("E6").Value = "  -2.35%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").Value = "'2.599.83"
$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("E13").Value = "  +3.33%  "

$ws.Range("D14").Value = "'3.043.61"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "'58.952.58"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").Value = "'20.43"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "'2.590.20"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'338.36"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").Value = "'10.09"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").Value = "'6.50"
$ws.Range("E22").Value = "  +2.73%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'66.07"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "'0.0₃0723"
$ws.Range("E30").Value = "  -3.19%  "

$ws.Range("D31").Value = "'5.95"
$ws.Range("E31").Value = "  -5.39%  "

$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").Value = "'18.68"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").Value = "'148.90"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("D37").Value = "'36.75"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "  -6.84%  "

$ws.Range("E41").Value = "  -0.62%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "'272.38"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").Value = "'10.75"
$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").Value = "'0.0517"
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("E48").Value = "  -1.75%  "

$ws.Range("D49").Value = "'1.970.76"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").Value = "'4.54"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("E51").Value = "  -0.38%  "
